$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New GPS survey row for Rockefeller Hall, appended after the existing
# Uris Library row (row 18).
$ws.Range("A19").Value = "Rockefeller Hall"
$ws.Range("B19").Value = 42.449079764950703
$ws.Range("C19").Value = -76.481934785842896
$ws.Range("D19").Value = 42.448739351674803
$ws.Range("E19").Value = -76.482149362564002
$ws.Range("F19").Value = 42.449459758980502
$ws.Range("G19").Value = -76.481623649597097
$ws.Range("H19").Formula = "=SQRT((F19-D19)^2+(G19-E19)^2)/2"
$ws.Range("I19").Value = 42.449428092899403
$ws.Range("J19").Value = -76.482160091400104
$ws.Range("K19").Value = 42.448755184891397
$ws.Range("L19").Value = -76.481537818908606
$ws.Range("M19").Formula = "=SQRT((K19-I19)^2+(L19-J19)^2)/2"
$ws.Range("N19").Formula = "=AVERAGE(H19,M19)"

# Leave the selection where the data-entry user ended up after filling
# the new row's formulas down from the previous row.
$ws.Range("C17:E18").Select()
